$d = $word.ActiveDocument

# Locate the specific "Checkin" paragraph that belongs to the calendar
# row for the 17th (the cell also contains 18 / 19 "Linked Lists (CH 4.1)"
# / 20 "Lab 4" "TextEdit Milestone 2 due end of day"). There are several
# "Checkin" paragraphs in the document, so we find the correct one by
# scanning the whole document's paragraphs collection (which indexes
# reliably) for the text "Checkin" that is immediately preceded by a
# paragraph containing just "17".
$targetIndex = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Checkin") {
        $prev = $d.Paragraphs.Item($i - 1)
        $prevTxt = $prev.Range.Text.TrimEnd([char]13, [char]7)
        if ($prevTxt -eq "17") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -eq -1) {
    Write-Output "ERROR: could not locate target paragraph"
} else {
    $p = $d.Paragraphs.Item($targetIndex)
    $searchRange = $d.Range($p.Range.Start, $p.Range.End)

    # Replace just this one occurrence of "Checkin" with "NO CLASS".
    $searchRange.Find.Execute("Checkin", $true, $false, $false, $false, $false, $true, 0, $false, "NO CLASS", 1)

    # Re-fetch the (now updated) paragraph range.
    $p = $d.Paragraphs.Item($targetIndex)

    # Move the "_GoBack" bookmark so that it sits right after the new
    # "NO CLASS" text instead of at the end of the document's last edited
    # location (after "TextEdit Milestone 2 due end of day").
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $bookmarkPos = $p.Range.End - 1
    $newBookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $newBookmarkRange)
}
